$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.070.46'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '3.050.95'
$ws.Range("E3").Value = '  -1.26%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").Value = '''582.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("D6").Value = '''150.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.50%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '''0.530'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.27%  '
$ws.Range("D9").Value = '3.052.28'
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("D10").Value = '''0.153'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.85%  '
$ws.Range("D11").Value = '''5.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.70%  '
$ws.Range("D12").Value = '''0.447'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.01%  '
$ws.Range("D13").Value = '''0.0000233'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.14%  '
$ws.Range("D14").Value = '''35.89'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.15%  '
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("D16").Value = '3.547.41'
$ws.Range("E16").Value = '  -1.50%  '
$ws.Range("D17").Value = '''7.11'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("D18").Value = '62.988.80'
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("D19").Value = '3.046.32'
$ws.Range("E19").Value = '  -1.22%  '
$ws.Range("D20").Value = '''477.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").Value = '''14.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.49%  '
$ws.Range("D22").Value = '''0.703'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.55%  '
$ws.Range("D23").Value = '''7.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").Value = '''2.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.97%  '
$ws.Range("D25").Value = '''81.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("D26").Value = '''12.58'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.44%  '
$ws.Range("D27").Value = '''10.57'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.19%  '
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("B29").Value = 'NEARProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D29").Value = '''7.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.97%  '
$ws.Range("B30").Value = 'FirstDigitalUSD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D30").Value = '''0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '''2.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''2.64'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.79%  '
$ws.Range("D33").Value = '''27.70'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.94%  '
$ws.Range("D34").Value = '''0.110'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.60%  '
$ws.Range("D35").Value = '''1.06'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.06%  '
$ws.Range("D36").Value = '0.0₃0807'
$ws.Range("E36").Value = '  -4.81%  '
$ws.Range("D37").Value = '''5.87'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.32%  '
$ws.Range("D38").Value = '''2.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.72%  '
$ws.Range("D39").Value = '''3.08'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.67%  '
$ws.Range("D40").Value = '''50.19'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.96%  '
$ws.Range("D41").Value = '''9.11'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.85%  '
$ws.Range("D42").Value = '''425.25'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.22%  '
$ws.Range("D43").Value = '''0.285'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("D44").Value = '''0.114'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.49%  '
$ws.Range("D45").Value = '2.830.81'
$ws.Range("E45").Value = '  +1.07%  '
$ws.Range("D46").Value = '''0.0360'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.64%  '
$ws.Range("D47").Value = '''38.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.95%  '
$ws.Range("D48").Value = '''127.78'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.67%  '
$ws.Range("D50").Value = '''24.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.99%  '
$ws.Range("E51").Value = '  -0.66%  '

Write-Output "done"